{"js": "// Update the date title and the 25 division-problem cells in the single\n// table, matching the diff:\n//   2024-09-09 Monday -> 2024-09-10 Tuesday\n//   plus 25 \"a\u00f7b=c, d\" cell rewrites (5 populated rows x 5 cols).\n\n// 1) Title paragraph (first paragraph in the body, outside the table).\nconst titleParagraphs = context.document.body.paragraphs;\ntitleParagraphs.load(\"items/text\");\nawait context.sync();\n\nconst titlePara = titleParagraphs.items[0];\nif (titlePara.text === \"2024-09-09 Monday\") {\n  titlePara.getRange().insertText(\"2024-09-10 Tuesday\", \"Replace\");\n}\n\n// 2) Table of division problems (first/only table in the document).\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\n// New values, same row/column layout as the existing table (rows with no\n// problems stay untouched -- we only overwrite the rows that already carry\n// text so blank spacer rows are left exactly as they are).\nconst newRows = {\n  0: [\"80\u00f76=13, 2\", \"21\u00f75=4, 1\", \"29\u00f79=3, 2\", \"23\u00f77=3, 2\", \"63\u00f79=7, 0\"],\n  4: [\"26\u00f74=6, 2\", \"58\u00f76=9, 4\", \"17\u00f74=4, 1\", \"65\u00f75=13, 0\", \"80\u00f72=40, 0\"],\n  8: [\"36\u00f75=7, 1\", \"97\u00f78=12, 1\", \"23\u00f78=2, 7\", \"72\u00f76=12, 0\", \"84\u00f76=14, 0\"],\n  12: [\"46\u00f74=11, 2\", \"60\u00f78=7, 4\", \"11\u00f78=1, 3\", \"17\u00f76=2, 5\", \"19\u00f78=2, 3\"],\n  16: [\"46\u00f78=5, 6\", \"75\u00f74=18, 3\", \"36\u00f79=4, 0\", \"59\u00f74=14, 3\", \"78\u00f74=19, 2\"],\n};\n\nfor (const rowIndexStr of Object.keys(newRows)) {\n  const rowIndex = Number(rowIndexStr);\n  const rowValues = newRows[rowIndex];\n  for (let col = 0; col < rowValues.length; col++) {\n    table.getCell(rowIndex, col).value = rowValues[col];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date title and the 25 division-problem cells in the single\n# table, matching the diff:\n#   2024-09-09 Monday -> 2024-09-10 Tuesday\n#   plus 25 \"a\u00f7b=c, d\" cell rewrites (5 populated rows x 5 cols).\n\n$d = $word.ActiveDocument\n\n# 1) Title paragraph (first paragraph in the body, outside the table).\n# Paragraph.Range.Text includes the trailing paragraph mark (chr 13), so\n# trim it before comparing; assigning the replacement text back leaves the\n# paragraph mark untouched.\n$titlePara = $d.Paragraphs.Item(1)\n$titleText = $titlePara.Range.Text.TrimEnd([char]13)\nif ($titleText -eq \"2024-09-09 Monday\") {\n    $titlePara.Range.Text = \"2024-09-10 Tuesday\"\n}\n\n# 2) Table of division problems (first/only table in the document).\n$table = $d.Tables.Item(1)\n\n# New values, same row/column layout as the existing table (rows with no\n# problems stay untouched -- only the rows that already carry text are\n# overwritten, so blank spacer rows are left exactly as they are). Word\n# table rows/columns are 1-based.\n$newRows = @{\n    1  = @(\"80\u00f76=13, 2\", \"21\u00f75=4, 1\", \"29\u00f79=3, 2\", \"23\u00f77=3, 2\", \"63\u00f79=7, 0\")\n    5  = @(\"26\u00f74=6, 2\", \"58\u00f76=9, 4\", \"17\u00f74=4, 1\", \"65\u00f75=13, 0\", \"80\u00f72=40, 0\")\n    9  = @(\"36\u00f75=7, 1\", \"97\u00f78=12, 1\", \"23\u00f78=2, 7\", \"72\u00f76=12, 0\", \"84\u00f76=14, 0\")\n    13 = @(\"46\u00f74=11, 2\", \"60\u00f78=7, 4\", \"11\u00f78=1, 3\", \"17\u00f76=2, 5\", \"19\u00f78=2, 3\")\n    17 = @(\"46\u00f78=5, 6\", \"75\u00f74=18, 3\", \"36\u00f79=4, 0\", \"59\u00f74=14, 3\", \"78\u00f74=19, 2\")\n}\n\nforeach ($rowIndex in $newRows.Keys) {\n    $rowValues = $newRows[$rowIndex]\n    for ($col = 1; $col -le $rowValues.Count; $col++) {\n        $table.Cell($rowIndex, $col).Range.Text = $rowValues[$col - 1]\n    }\n}\n"}
